# Nalco aluminium-ingot price sheet: prepend the new top row (13-01-2026),
# pushing every existing row down by one, and extend the table with one more
# trailing row duplicating the previous last row (the pre-shift last row ends
# up twice: the oldest still-current row and the newly appended historical
# copy).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 160
$newLastRow = $lastRow + 1

# Columns A (Date) and E (Circular Date) hold dd-mm-yyyy text that Excel's
# COM layer will happily misparse as a US-style m-d-y date the moment it is
# ambiguous (day <= 12). Force those columns to Text format up front so every
# write below is stored as the literal string, matching the original
# inlineStr cells.
$ws.Range("A2:A$newLastRow").NumberFormat = "@"
$ws.Range("E2:E$newLastRow").NumberFormat = "@"

# Shift existing data rows 2..160 down to 3..161, working bottom-up so we
# never clobber a row before it has been read.
for ($r = $lastRow; $r -ge 2; $r--) {
    $dest = $r + 1
    for ($c = 1; $c -le 6; $c++) {
        $src = $ws.Cells.Item($r, $c)
        $ws.Cells.Item($dest, $c).Value = $src.Value()
    }
}

# New row 2: latest price update, carried over from the (now shifted-down)
# former row 2 values, with just the date bumped to 13-01-2026.
$ws.Cells.Item(2, 1).Value = "13-01-2026"
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 307.25
$ws.Cells.Item(2, 5).Value = "01-01-2026"
$ws.Cells.Item(2, 6).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-01-2026.pdf"

# Row 161 is brand new (previously blank). Give its text/number columns the
# same centered alignment used throughout the table so they pick up the
# matching shared cell styles, then fill in the duplicate of the old last row.
$ws.Cells.Item($newLastRow, 2).HorizontalAlignment = -4108
$ws.Cells.Item($newLastRow, 2).VerticalAlignment = -4108
$ws.Cells.Item($newLastRow, 3).HorizontalAlignment = -4108
$ws.Cells.Item($newLastRow, 3).VerticalAlignment = -4108
$ws.Cells.Item($newLastRow, 4).NumberFormat = "0.000"
$ws.Cells.Item($newLastRow, 4).HorizontalAlignment = -4108
$ws.Cells.Item($newLastRow, 4).VerticalAlignment = -4108
$ws.Cells.Item($newLastRow, 1).HorizontalAlignment = -4108
$ws.Cells.Item($newLastRow, 1).VerticalAlignment = -4108
$ws.Cells.Item($newLastRow, 5).HorizontalAlignment = -4108
$ws.Cells.Item($newLastRow, 5).VerticalAlignment = -4108

$ws.Cells.Item($newLastRow, 1).Value = "07-08-2025"
$ws.Cells.Item($newLastRow, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item($newLastRow, 3).Value = "IE07"
$ws.Cells.Item($newLastRow, 4).Value = 268.25
$ws.Cells.Item($newLastRow, 5).Value = "07-08-2025"
$ws.Cells.Item($newLastRow, 6).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# The Circular Link column (F) is backed by real Hyperlink objects whose
# target does NOT follow a plain .Value write, so rebuild them from scratch:
# clear every hyperlink on the sheet (cell-scoped Delete() clears them all in
# this engine) and re-add one per row, now that every F-cell holds its final
# URL text.
$ws.Cells.Item(2, 6).Hyperlinks.Delete()
for ($r = 2; $r -le $newLastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $cell.Value())
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

Write-Output ("rows=" + $ws.Cells.Item($newLastRow, 1).Value() + " hyperlinks=" + $ws.Hyperlinks.Count)
